$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Introduction paragraph: fix the typo "leviador" -> "levitador"
#    (the user inserted a missing "t"), and merge the paragraph near
#    "cada componente..." where a stray run split (and the old
#    "_GoBack" bookmark) used to sit.
# ------------------------------------------------------------------

# Merge "No limitarse a u" / "na descripcion..." back into a single
# run - this also drops the old "_GoBack" bookmark that used to live
# at that split point, since the edit moved on to a new location.
$d.Content.Find.Execute("No limitarse a una descripción básica", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "No limitarse a una descripción básica", 2)

# Fix "leviador" -> "levitador".
$d.Content.Find.Execute("Un leviador neumático es", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Un levitador neumático es", 2)

# ------------------------------------------------------------------
# 2) Drop the "_GoBack" bookmark at its new (last-edited) location:
#    right after "...controlador P" (mid-word in "PID"), which is
#    where Word leaves it after the user's most recent edit.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("controlador P", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
